$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

$metaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our Don Juan's Peppers review and play for free. Enjoy a Mexican-themed slot game with big win opportunities and a generous free spins mode.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$insertPoint.InsertXML($metaXml)

# The inserted fragment carries a trailing placeholder paragraph so the new
# "Meta description" paragraph doesn't swallow the following heading; remove
# that now-empty placeholder paragraph.
$placeholder = $d.Paragraphs(3)
$placeholder.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Near the end of the document: remove the bold "Play Don Juan's Peppers
#    free online slot game" paragraph and rewrite the remaining italic
#    paragraph's text into the new AI-image-prompt text.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPromoPara = $d.Paragraphs($count - 1)
$boldPromoPara.Range.Delete()

$count2 = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($count2)
$italicRange = $italicPara.Range
$italicTextRange = $d.Range($italicRange.Start, $italicRange.End - 1)
$italicTextRange.Text = "Please create a feature image fitting the game `"Don Juan's Peppers`". The image should be in a cartoon style and feature a happy Maya warrior with glasses."
